$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-04-04 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-05 Wednesday", 2) | Out-Null

# Update the multiplication table cells by position to avoid any
# ambiguity from values that coincide with other cells old/new text
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "96×31=2976"
$t.Cell(1, 2).Range.Text = "76×39=2964"
$t.Cell(1, 3).Range.Text = "88×22=1936"
$t.Cell(1, 4).Range.Text = "66×30=1980"
$t.Cell(1, 5).Range.Text = "71×18=1278"

$t.Cell(2, 1).Range.Text = "81×53=4293"
$t.Cell(2, 2).Range.Text = "47×46=2162"
$t.Cell(2, 3).Range.Text = "42×48=2016"
$t.Cell(2, 4).Range.Text = "95×35=3325"
$t.Cell(2, 5).Range.Text = "21×88=1848"

$t.Cell(3, 1).Range.Text = "96×23=2208"
$t.Cell(3, 2).Range.Text = "40×63=2520"
$t.Cell(3, 3).Range.Text = "86×13=1118"
$t.Cell(3, 4).Range.Text = "51×67=3417"
$t.Cell(3, 5).Range.Text = "66×82=5412"

$t.Cell(4, 1).Range.Text = "48×49=2352"
$t.Cell(4, 2).Range.Text = "28×59=1652"
$t.Cell(4, 3).Range.Text = "55×43=2365"
$t.Cell(4, 4).Range.Text = "41×32=1312"
$t.Cell(4, 5).Range.Text = "42×20=840"

$t.Cell(5, 1).Range.Text = "63×52=3276"
$t.Cell(5, 2).Range.Text = "30×44=1320"
$t.Cell(5, 3).Range.Text = "16×58=928"
$t.Cell(5, 4).Range.Text = "18×56=1008"
$t.Cell(5, 5).Range.Text = "23×61=1403"

$t.Cell(6, 1).Range.Text = "67×45=3015"
$t.Cell(6, 2).Range.Text = "85×10=850"
$t.Cell(6, 3).Range.Text = "51×67=3417"
$t.Cell(6, 4).Range.Text = "96×28=2688"
$t.Cell(6, 5).Range.Text = "45×56=2520"

$t.Cell(7, 1).Range.Text = "13×26=338"
$t.Cell(7, 2).Range.Text = "32×24=768"
$t.Cell(7, 3).Range.Text = "46×48=2208"
$t.Cell(7, 4).Range.Text = "92×97=8924"
$t.Cell(7, 5).Range.Text = "58×15=870"

$t.Cell(8, 1).Range.Text = "96×27=2592"
$t.Cell(8, 2).Range.Text = "73×38=2774"
$t.Cell(8, 3).Range.Text = "68×79=5372"
$t.Cell(8, 4).Range.Text = "51×64=3264"
$t.Cell(8, 5).Range.Text = "50×60=3000"

$t.Cell(9, 1).Range.Text = "87×79=6873"
$t.Cell(9, 2).Range.Text = "46×30=1380"
$t.Cell(9, 3).Range.Text = "37×47=1739"
$t.Cell(9, 4).Range.Text = "79×12=948"
$t.Cell(9, 5).Range.Text = "69×32=2208"

$t.Cell(10, 1).Range.Text = "41×55=2255"
$t.Cell(10, 2).Range.Text = "19×16=304"
$t.Cell(10, 3).Range.Text = "54×34=1836"
$t.Cell(10, 4).Range.Text = "53×36=1908"
$t.Cell(10, 5).Range.Text = "96×64=6144"

$t.Cell(11, 1).Range.Text = "47×17=799"
$t.Cell(11, 2).Range.Text = "44×68=2992"
$t.Cell(11, 3).Range.Text = "81×27=2187"
$t.Cell(11, 4).Range.Text = "81×30=2430"
$t.Cell(11, 5).Range.Text = "79×22=1738"

$t.Cell(12, 1).Range.Text = "13×12=156"
$t.Cell(12, 2).Range.Text = "24×21=504"
$t.Cell(12, 3).Range.Text = "90×74=6660"
$t.Cell(12, 4).Range.Text = "59×73=4307"
$t.Cell(12, 5).Range.Text = "26×61=1586"

$t.Cell(13, 1).Range.Text = "90×20=1800"
$t.Cell(13, 2).Range.Text = "46×52=2392"
$t.Cell(13, 3).Range.Text = "54×18=972"
$t.Cell(13, 4).Range.Text = "90×25=2250"
$t.Cell(13, 5).Range.Text = "95×77=7315"

$t.Cell(14, 1).Range.Text = "65×96=6240"
$t.Cell(14, 2).Range.Text = "56×97=5432"
$t.Cell(14, 3).Range.Text = "41×54=2214"
$t.Cell(14, 4).Range.Text = "56×81=4536"
$t.Cell(14, 5).Range.Text = "31×58=1798"

$t.Cell(15, 1).Range.Text = "80×25=2000"
$t.Cell(15, 2).Range.Text = "59×23=1357"
$t.Cell(15, 3).Range.Text = "22×83=1826"
$t.Cell(15, 4).Range.Text = "30×31=930"
$t.Cell(15, 5).Range.Text = "88×65=5720"

$t.Cell(16, 1).Range.Text = "11×69=759"
$t.Cell(16, 2).Range.Text = "92×97=8924"
$t.Cell(16, 3).Range.Text = "90×61=5490"
$t.Cell(16, 4).Range.Text = "85×85=7225"
$t.Cell(16, 5).Range.Text = "22×29=638"

$t.Cell(17, 1).Range.Text = "49×13=637"
$t.Cell(17, 2).Range.Text = "64×11=704"
$t.Cell(17, 3).Range.Text = "75×54=4050"
$t.Cell(17, 4).Range.Text = "49×12=588"
$t.Cell(17, 5).Range.Text = "50×44=2200"

$t.Cell(18, 1).Range.Text = "39×75=2925"
$t.Cell(18, 2).Range.Text = "78×29=2262"
$t.Cell(18, 3).Range.Text = "45×95=4275"
$t.Cell(18, 4).Range.Text = "13×73=949"
$t.Cell(18, 5).Range.Text = "58×38=2204"

$t.Cell(19, 1).Range.Text = "17×50=850"
$t.Cell(19, 2).Range.Text = "88×81=7128"
$t.Cell(19, 3).Range.Text = "59×96=5664"
$t.Cell(19, 4).Range.Text = "20×76=1520"
$t.Cell(19, 5).Range.Text = "58×39=2262"

$t.Cell(20, 1).Range.Text = "96×92=8832"
$t.Cell(20, 2).Range.Text = "56×23=1288"
$t.Cell(20, 3).Range.Text = "98×68=6664"
$t.Cell(20, 4).Range.Text = "34×76=2584"
$t.Cell(20, 5).Range.Text = "76×49=3724"
